$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E (price / volume%) stay as literal text so values
# like "1.000", "0.9997", "29.127.75" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.127.75"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "1.832.08"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "241.24"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").Value = "0.6634"
$ws.Range("E6").Value = "  -2.68%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "0.07417"
$ws.Range("E8").Value = "  -0.40%  "

$ws.Range("D9").Value = "0.2935"
$ws.Range("E9").Value = "  -1.91%  "

$ws.Range("D10").Value = "22.68"
$ws.Range("E10").Value = "  -2.32%  "

$ws.Range("D11").Value = "0.07733"
$ws.Range("E11").Value = "  +1.13%  "

$ws.Range("D12").Value = "1.844.02"
$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("D13").Value = "4.988"
$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").Value = "0.6679"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").Value = "82.77"
$ws.Range("E15").Value = "  -5.36%  "

$ws.Range("D16").Value = "6.099"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").Value = "0.000008348"
$ws.Range("E17").Value = "  +1.69%  "

$ws.Range("D18").Value = "29.148.53"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").Value = "227.26"
$ws.Range("E19").Value = "  -1.31%  "

$ws.Range("D20").Value = "12.47"
$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "7.153"
$ws.Range("E22").Value = "  -2.59%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "159.65"
$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("D25").Value = "8.614"
$ws.Range("E25").Value = "  -1.05%  "

$ws.Range("D26").Value = "0.1401"
$ws.Range("E26").Value = "  -1.93%  "

$ws.Range("D27").Value = "17.95"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").Value = "1.509"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").Value = "4.111"
$ws.Range("E29").Value = "  -3.30%  "

$ws.Range("D30").Value = "4.042"
$ws.Range("E30").Value = "  -2.32%  "

$ws.Range("D31").Value = "1.195"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D32").Value = "0.05312"
$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("D33").Value = "1.869"
$ws.Range("E33").Value = "  +1.21%  "

$ws.Range("D34").Value = "0.7515"
$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("D35").Value = "1.136"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").Value = "2.637"
$ws.Range("E36").Value = "  -1.70%  "

$ws.Range("D37").Value = "1.281.38"
$ws.Range("E37").Value = "  -2.23%  "

$ws.Range("D38").Value = "0.01796"
$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("D39").Value = "2.735"
$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("D40").Value = "0.9285"
$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("D41").Value = "0.08739"
$ws.Range("E41").Value = "  +9.11%  "

$ws.Range("D42").Value = "5.928"
$ws.Range("E42").Value = "  -2.11%  "

$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "102.10"
$ws.Range("E44").Value = "  -2.59%  "

$ws.Range("D45").Value = "1.983.77"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").Value = "0.5142"
$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("D47").Value = "1.765"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").Value = "63.22"
$ws.Range("E49").Value = "  -1.32%  "

$ws.Range("D50").Value = "0.05896"
$ws.Range("E50").Value = "  -0.91%  "

# Row 51: coin replaced (EnergySwap -> Aptos) along with its link/price/volume
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "6.784"
$ws.Range("E51").Value = "  -1.50%  "
